# Grades workbook update: add HW 7-11 (CPP) grade groups for the single
# student row, mirroring the existing "Assignment / Grade / Comments"
# column groups (pattern already used for groups 3-7, e.g. columns J:M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an existing fully-formatted group (J1:M2 = Assignment/Grade/Comments
# + trailing separator column, both header & data row formatting) as the
# template and copy it into the five new slots so every new cell picks up
# the correct style (centered header s=1 / bordered data s=4 / separator
# s=2 & s=5) instead of the bare default column style.
$template = $ws.Range("J1:M2")
$newGroupStarts = @("AC1", "AG1", "AK1", "AO1", "AS1")
foreach ($startCell in $newGroupStarts) {
    $template.Copy($ws.Range($startCell)) | Out-Null
}

# Fill in the real values for the five new assignment groups.
$ws.Range("AC1").Value2 = "Assignment"
$ws.Range("AD1").Value2 = "Grade"
$ws.Range("AE1").Value2 = "Comments"

$ws.Range("AC2").Value2 = "7 CPP"
$ws.Range("AD2").Value2 = 80
$ws.Range("AE2").Value2 = "comments in the code"

$ws.Range("AG1").Value2 = "Assignment"
$ws.Range("AH1").Value2 = "Grade"
$ws.Range("AI1").Value2 = "Comments"

$ws.Range("AG2").Value2 = "8 CPP"
$ws.Range("AH2").Value2 = 85
$ws.Range("AI2").Value2 = "comments in the code"

$ws.Range("AK1").Value2 = "Assignment"
$ws.Range("AL1").Value2 = "Grade"
$ws.Range("AM1").Value2 = "Comments"

$ws.Range("AK2").Value2 = "9 CPP"
$ws.Range("AL2").Value2 = 100
$ws.Range("AM2").Value2 = "excellent"

$ws.Range("AO1").Value2 = "Assignment"
$ws.Range("AP1").Value2 = "Grade"
$ws.Range("AQ1").Value2 = "Comments"

$ws.Range("AO2").Value2 = "10 CPP"
$ws.Range("AP2").Value2 = 0
$ws.Range("AQ2").Value2 = "you need to submit this separately"

$ws.Range("AS1").Value2 = "Assignment"
$ws.Range("AT1").Value2 = "Grade"
$ws.Range("AU1").Value2 = "Comments"

$ws.Range("AS2").Value2 = "11 CPP"
$ws.Range("AT2").Value2 = 95
$ws.Range("AU2").Value2 = "very good, see my example in class 11 about the virtual Clone() method"

# Header row grew taller (wrapped "Assignment/Grade/Comments" labels) and
# the data row grew to fit the longer comments.
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 150

# Selection ends on the last newly-filled cell.
$ws.Range("AV2").Select() | Out-Null
